$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the season-record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting used by the rest of row 1 (bold, bordered,
# centered/top-aligned) by copying formats from an existing header cell.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Populate the season record (Wins/Losses/Ties) for every data row.
for ($row = 2; $row -le 37; $row++) {
    $ws.Cells.Item($row, 30).Value = 96
    $ws.Cells.Item($row, 31).Value = 67
    $ws.Cells.Item($row, 32).Value = 0
}
